$d = $word.ActiveDocument

# 1. Replace the long paragraph asking about speaking with Christian Hamm,
#    James Faulkner and Jake Downie with a short closing sentence.
#    Find/Execute replaces the whole matched span (which is split across
#    many runs) with a single new run, instead of leaving the old runs
#    behind.
$oldText = "Would you also like to speak with Christian Hamm, James Faulkner and Jake" `
    + " Downie at SiteMax Systems? They are amazing people, and they always" `
    + " looking to improve their sales and help the company grow, and during" `
    + " COVID-19 a help like this may be beneficial for the company.  "
$newText = "I hope to speak with you soon."

$null = $d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, `
    $true, 1, $false, $newText, 2)

# 2. That paragraph used to be followed by two blank paragraphs; the edit
#    collapses them down to one, so remove the first of the pair.
$d.Paragraphs.Item(8).Range.Delete()

# 3. Add a new blank paragraph right after the closing "Hyungmo Gu" line,
#    at the very end of the document body.
$d.Paragraphs.Last.Range.InsertParagraphAfter()
